$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Range("C9").Interior.ThemeColor = 4
$ws.Range("C10").Interior.ThemeColor = 5
$ws.Range("C11").Interior.ThemeColor = 6
$ws.Range("C12").Interior.ThemeColor = 7
$ws.Range("C13").Interior.ThemeColor = 8
$ws.Range("C14").Interior.ThemeColor = 9
$ws.Range("C15").Interior.ThemeColor = 10
$ws.Range("C16").Interior.ThemeColor = 11
$ws.Range("C17").Interior.ThemeColor = 12
Write-Host "done"
